$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.424719
$ws.Range("N2").Value = 4.274157
$ws.Range("O2").Value = 0.07423298812267187
$ws.Range("P2").Value = 0.07423298812267187
$ws.Range("Q2").Value = 0.1035637739226667
$ws.Range("R2").Value = 0.9320739653039999
$ws.Range("S2").Value = 0.07423298812267187
$ws.Range("T2").Value = 0.07423298812267187

# Row 3 updates
$ws.Range("O3").Value = 0.5596266124066729
$ws.Range("P3").Value = 0.5596266124066729
$ws.Range("S3").Value = 0.5596266124066729
$ws.Range("T3").Value = 0.5596266124066729

# Row 4 updates
$ws.Range("M4").Value = 7.027161333333335
$ws.Range("O4").Value = 0.3661403994706553
$ws.Range("P4").Value = 0.3661403994706552
$ws.Range("Q4").Value = 0.5108090420942223
$ws.Range("R4").Value = 4.597281378848001
$ws.Range("S4").Value = 0.3661403994706553
$ws.Range("T4").Value = 0.3661403994706552
